# Journal de bord - add new entry for "Création de la maquette" and
# move the active selection, mirroring the author's edits captured in the
# commit "Update Journal de bord, Journal de travail + Add Maquette.png".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 was an empty, pre-formatted template row. The author filled it in
# with a new journal entry. Copy the formatting from the row above (row 5,
# the last populated data row) so the new row matches the data-row style
# (borders/fill) instead of keeping the generic empty-template style.
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new journal entry values.
$ws.Range("A6").Value = "Gabriel Pereira"
$ws.Range("B6").Value = 43896
$ws.Range("C6").Value = 0.99444444444444446
$ws.Range("D6").Value = "Bataille Navale"

# Set F6 (Description) before E6 (Tache) so the new shared strings are
# appended in the same order as the reference workbook:
#   index 16 -> "Création de la maquette" (Description)
#   index 17 -> "Création de maquette"    (Tache)
$ws.Range("F6").Value = "Création de la maquette"
$ws.Range("E6").Value = "Création de maquette"

# Move the active selection to C9, matching the saved selection state.
$ws.Range("C9").Select()
